$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Цена за все экраны" -> "Цена за все товары" (column E header) ---
$ws.Range("E1").Value = "Цена за все товары"

# --- Row 2 (candle): quantity 2 -> 3, recompute line total 500*3 = 1500 ---
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 1500

# --- Row 3: new product "sylvia-breitenberg", price 615, qty 5, line total 3075 ---
$ws.Range("A3").Value = "sylvia-breitenberg"
$ws.Range("B3").Value = 615
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 3075

# The old grand-total cell (E3 = 1000) moves down to E4 with the new total
$ws.Range("E3").ClearContents()
$ws.Range("E4").Value = 4575

# --- Column widths: A 16->28, C 28->40 ---
# ColumnWidth is in character units; the stored <col width> attribute adds the
# standard ~5px padding (5/6 of a character at the default font), so subtract
# that offset here to land on the exact target stored width.
$ws.Columns.Item(1).ColumnWidth = 28 - 5/6
$ws.Columns.Item(3).ColumnWidth = 40 - 5/6

Write-Host "done"
